$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update participated count for "III year cse" department
$ws.Range("B2").Value = 128

# Update the computed percentage (participated/total*100) to match
$ws.Range("D2").Value = 70.72
